# Applies the shared-strings / cell value changes described by the diff:
# a new entry "1033242 - Fábio Herbst Florenzano" is introduced and the
# existing table/entry shuffle ends up changing the text shown in a handful
# of B/C value cells on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (B13/C13) becomes the text that used to live in B15/C15 ("Semestral").
# Copy it first, before B15/C15 themselves get overwritten below.
$ws.Range("B15").Copy()
$ws.Range("B13").PasteSpecial(-4104)
$ws.Range("C15").Copy()
$ws.Range("C13").PasteSpecial(-4104)

# Row 15 (B15/C15) becomes the text that lives in B8/C8 ("01/01/2023"), copied
# as a plain text cell so Excel doesn't reinterpret it as a date value.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4104)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4104)

# Rows 10 and 18 (B/C) get the newly introduced instructor string.
$ws.Range("B10").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C10").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C18").Value = "1033242 - Fábio Herbst Florenzano"

$excel.CutCopyMode = $false
